$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to text format so numeric-looking strings
# (e.g. "243.96", "0.620") keep their original textual representation
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "36.370.93"
$ws.Range("E2").Value = "  -2.82%  "
$ws.Range("D3").Value = "1.964.92"
$ws.Range("E3").Value = "  -4.43%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "243.96"
$ws.Range("E5").Value = "  -3.40%  "
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("D7").Value = "58.20"
$ws.Range("E7").Value = "  -11.82%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -6.89%  "
$ws.Range("D10").Value = "55.84"
$ws.Range("E10").Value = "  -5.87%  "
$ws.Range("D11").Value = "0.0852"
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "22.26"
$ws.Range("E13").Value = "  -6.35%  "
$ws.Range("E14").Value = "  -9.26%  "
$ws.Range("D15").Value = "2.250.79"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("D16").Value = "13.50"
$ws.Range("E16").Value = "  -8.90%  "
$ws.Range("D17").Value = "5.36"
$ws.Range("E17").Value = "  -5.88%  "
$ws.Range("D18").Value = "1.981.44"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "36.312.36"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").Value = "71.27"
$ws.Range("E20").Value = "  -2.80%  "
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "5.15"
$ws.Range("E22").Value = "  -6.58%  "
$ws.Range("D23").Value = "231.30"
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "2.54"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("E26").Value = "  -4.77%  "
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").Value = "  -4.58%  "
$ws.Range("D28").Value = "166.93"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("D29").Value = "20.00"
$ws.Range("E29").Value = "  -4.13%  "
$ws.Range("D30").Value = "0.124"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("D33").Value = "4.76"
$ws.Range("E33").Value = "  -8.96%  "
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("D35").Value = "4.36"
$ws.Range("E35").Value = "  -6.71%  "
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").Value = "5.97"
$ws.Range("E38").Value = "  -6.90%  "
$ws.Range("D39").Value = "2.16"
$ws.Range("E39").Value = "  -10.69%  "
$ws.Range("D40").Value = "2.94"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").Value = "0.0963"
$ws.Range("E41").Value = "  -5.70%  "
$ws.Range("E42").Value = "  -5.16%  "
$ws.Range("E43").Value = "  -8.39%  "
$ws.Range("D44").Value = "0.0211"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("E45").Value = "  -9.65%  "
$ws.Range("D46").Value = "15.77"
$ws.Range("E46").Value = "  -9.43%  "
$ws.Range("D47").Value = "88.91"
$ws.Range("E47").Value = "  -7.13%  "
$ws.Range("D48").Value = "1.349.08"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").Value = "7.29"
$ws.Range("E49").Value = "  -7.68%  "
$ws.Range("E50").Value = "  -4.32%  "
$ws.Range("D51").Value = "45.13"
$ws.Range("E51").Value = "  -3.57%  "

# Reset column D style back to default (Normal) so no stray formatting
# is left behind, matching the original workbook styling
$ws.Range("D2:D51").Style = "Normal"
